$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 395.30768
$ws.Range("I33").Value = 395.30768
$ws.Range("K33").Value = 395.30768
$ws.Range("M33").Value = -166.30768
# Row 62
$ws.Range("H62").Value = 14303789
$ws.Range("I62").Value = 17878236
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 17878236
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -17877612
$ws.Range("N62").Value = -7248
# Row 65
$ws.Range("H65").Value = 14303789
$ws.Range("I65").Value = 17878236
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 89391180
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -89388060
$ws.Range("N65").Value = -36240
# Row 98
$ws.Range("H98").Value = 2568.9644
$ws.Range("I98").Value = 2330.5417
$ws.Range("K98").Value = 2330.5417
$ws.Range("M98").Value = -832.5417000000002
# Row 107
$ws.Range("H107").Value = 1877.8148
$ws.Range("I107").Value = 930.34784
$ws.Range("J107").Value = 7325.75
$ws.Range("K107").Value = 930.34784
$ws.Range("L107").Value = 7325.75
$ws.Range("M107").Value = 989.65216
$ws.Range("N107").Value = -11165.75
# Row 122
$ws.Range("H122").Value = 2568.9644
$ws.Range("I122").Value = 2330.5417
$ws.Range("K122").Value = 6991.625100000001
$ws.Range("M122").Value = -4541.625100000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 15886190
$ws.Range("I80").Value = 1425.1
$ws.Range("J80").Value = 30326884
$ws.Range("K80").Value = 1425.1
$ws.Range("L80").Value = 30326884
$ws.Range("M80").Value = -427.0999999999999
$ws.Range("N80").Value = -30328880
# Row 83
$ws.Range("H83").Value = 15886190
$ws.Range("I83").Value = 1425.1
$ws.Range("J83").Value = 30326884
$ws.Range("K83").Value = 7125.5
$ws.Range("L83").Value = 151634420
$ws.Range("M83").Value = -2133.5
$ws.Range("N83").Value = -151644404
# Row 99
$ws.Range("H99").Value = 6865.0527
$ws.Range("I99").Value = 7155.1953
$ws.Range("J99").Value = 6121.5625
$ws.Range("K99").Value = 7155.1953
$ws.Range("L99").Value = 6121.5625
$ws.Range("M99").Value = -5657.1953
$ws.Range("N99").Value = -9117.5625

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 94
$ws.Range("H94").Value = 45458476
$ws.Range("I94").Value = 71430640
$ws.Range("J94").Value = 7180
$ws.Range("K94").Value = 71430640
$ws.Range("L94").Value = 7180
$ws.Range("M94").Value = -71430189
$ws.Range("N94").Value = -8082
# Row 107
$ws.Range("H107").Value = 982.8333
$ws.Range("I107").Value = 899.6
$ws.Range("K107").Value = 899.6
$ws.Range("M107").Value = 1020.4
# Row 122
$ws.Range("H122").Value = 3041.8667
$ws.Range("I122").Value = 2538
$ws.Range("J122").Value = 4049.6
$ws.Range("K122").Value = 7614
$ws.Range("L122").Value = 12148.8
$ws.Range("M122").Value = -5164
$ws.Range("N122").Value = -17048.8
# Row 132
$ws.Range("H132").Value = 4270.3
$ws.Range("I132").Value = 3337.0222
$ws.Range("K132").Value = 10011.0666
$ws.Range("M132").Value = -7481.0666
# Row 134
$ws.Range("H134").Value = 71440184
$ws.Range("I134").Value = 125009220
$ws.Range("J134").Value = 14815.833
$ws.Range("K134").Value = 375027660
$ws.Range("L134").Value = 44447.499
$ws.Range("M134").Value = -375025125
$ws.Range("N134").Value = -49517.499

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 1332
$ws.Range("I8").Value = 1332
$ws.Range("K8").Value = 3996
$ws.Range("M8").Value = -3857
# Row 22
$ws.Range("H22").Value = 2150
$ws.Range("I22").Value = 599
$ws.Range("J22").Value = 2667
$ws.Range("K22").Value = 1797
$ws.Range("L22").Value = 8001
$ws.Range("M22").Value = -1628
$ws.Range("N22").Value = -8339
# Row 27
$ws.Range("H27").Value = 2150
$ws.Range("I27").Value = 599
$ws.Range("J27").Value = 2667
$ws.Range("K27").Value = 1797
$ws.Range("L27").Value = 8001
$ws.Range("M27").Value = -1695
$ws.Range("N27").Value = -8205
# Row 107
$ws.Range("H107").Value = 640.7895
$ws.Range("J107").Value = 670.13336
$ws.Range("L107").Value = 2010.40008
$ws.Range("N107").Value = -5850.40008
# Row 110
$ws.Range("H110").Value = 22562.334
$ws.Range("I110").Value = 11587
$ws.Range("K110").Value = 34761
$ws.Range("M110").Value = -30671
# Row 114
$ws.Range("H114").Value = 111889
$ws.Range("J114").Value = 333638
$ws.Range("L114").Value = 1000914
$ws.Range("N114").Value = -1007422
# Row 117
$ws.Range("H117").Value = 1666.0667
$ws.Range("J117").Value = 1229
$ws.Range("L117").Value = 3687
$ws.Range("N117").Value = -10571
# Row 129
$ws.Range("H129").Value = 9805764
$ws.Range("J129").Value = 15153952
$ws.Range("L129").Value = 45461856
$ws.Range("N129").Value = -45471856
# Row 131
$ws.Range("H131").Value = 24763920
$ws.Range("J131").Value = 19611224
$ws.Range("L131").Value = 58833672
$ws.Range("N131").Value = -58843752
# Row 139
$ws.Range("H139").Value = 21741464
$ws.Range("J139").Value = 3391.6667
$ws.Range("L139").Value = 10175.0001
$ws.Range("N139").Value = -20455.0001
# Row 140
$ws.Range("H140").Value = 38692116
$ws.Range("I140").Value = 45140070
$ws.Range("J140").Value = 4371.25
$ws.Range("K140").Value = 135420210
$ws.Range("L140").Value = 13113.75
$ws.Range("M140").Value = -135415030
$ws.Range("N140").Value = -23473.75

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 2348.7317
$ws.Range("I102").Value = 1848.3914
$ws.Range("J102").Value = 2988.0557
$ws.Range("K102").Value = 1848.3914
$ws.Range("L102").Value = 2988.0557
$ws.Range("M102").Value = -226.3914
$ws.Range("N102").Value = -6232.0557
# Row 122
$ws.Range("H122").Value = 3881.262
$ws.Range("I122").Value = 2912.0645
$ws.Range("J122").Value = 6612.636
$ws.Range("K122").Value = 8736.193499999999
$ws.Range("L122").Value = 19837.908
$ws.Range("M122").Value = -6286.193499999999
$ws.Range("N122").Value = -24737.908

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 6002.3
$ws.Range("J40").Value = 7887.2
$ws.Range("L40").Value = 7887.2
$ws.Range("N40").Value = -8159.2
# Row 43
$ws.Range("H43").Value = 2510624.5
$ws.Range("J43").Value = 2510624.5
$ws.Range("L43").Value = 2510624.5
$ws.Range("N43").Value = -2511010.5
# Row 82
$ws.Range("H82").Value = 4220
$ws.Range("I82").Value = 2400
$ws.Range("J82").Value = 6040
$ws.Range("K82").Value = 2400
$ws.Range("L82").Value = 6040
$ws.Range("M82").Value = -2039
$ws.Range("N82").Value = -6762
# Row 85
$ws.Range("H85").Value = 4220
$ws.Range("I85").Value = 2400
$ws.Range("J85").Value = 6040
$ws.Range("K85").Value = 2400
$ws.Range("L85").Value = 6040
$ws.Range("M85").Value = -1152
$ws.Range("N85").Value = -8536
# Row 100
$ws.Range("H100").Value = 3837.3076
$ws.Range("I100").Value = 7750
$ws.Range("J100").Value = 2098.3333
$ws.Range("K100").Value = 7750
$ws.Range("L100").Value = 2098.3333
$ws.Range("M100").Value = -7209
$ws.Range("N100").Value = -3180.3333
# Row 122
$ws.Range("H122").Value = 2697.8064
$ws.Range("I122").Value = 2545.6956
$ws.Range("K122").Value = 7637.0868
$ws.Range("M122").Value = -5187.0868

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 119
$ws.Range("H119").Value = 66999.8
$ws.Range("J119").Value = 66999.8
$ws.Range("L119").Value = 66999.8
$ws.Range("N119").Value = -76675.8
# Row 132
$ws.Range("H132").Value = 5684.8716
$ws.Range("I132").Value = 3049.743
$ws.Range("J132").Value = 28742.25
$ws.Range("K132").Value = 9149.228999999999
$ws.Range("L132").Value = 86226.75
$ws.Range("M132").Value = -6619.228999999999
$ws.Range("N132").Value = -91286.75
# Row 136
$ws.Range("H136").Value = 13894976
$ws.Range("I136").Value = 20834672
$ws.Range("J136").Value = 15585.417
$ws.Range("K136").Value = 62504016
$ws.Range("L136").Value = 46756.251
$ws.Range("M136").Value = -62501466
$ws.Range("N136").Value = -51856.251
